# Adds a new column V to the deceased_cases worksheet containing the
# 07-10-2020 daily figures, mirroring the existing U (06-10-2020) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: same style as U1, with the new date label stored as text
# (not auto-converted to a date serial number).
$ws.Range("V1").Style = $ws.Range("U1").Style
$ws.Range("V1").NumberFormat = "@"
$ws.Range("V1").Value = "07-10-2020"

# Per-state deceased counts for 07-10-2020, keyed by row number.
$values = @{
    2  = 54
    3  = 6052
    4  = 20
    5  = 778
    6  = 925
    7  = 180
    8  = 1104
    9  = 2
    10 = 5581
    11 = 468
    12 = 3519
    13 = 1509
    14 = 229
    15 = 1268
    16 = 757
    17 = 9461
    18 = 884
    19 = 61
    20 = 2488
    21 = 38717
    22 = 78
    23 = 60
    24 = 0
    25 = 17
    26 = 940
    27 = 546
    28 = 3679
    29 = 1574
    30 = 49
    31 = 9917
    32 = 1189
    33 = 301
    34 = 677
    35 = 6153
    36 = 5318
}

foreach ($row in 2..36) {
    $ws.Cells.Item($row, 22).Value = $values[$row]
}
